$d = $word.ActiveDocument
$r = $d.Content
$r.Find.Execute("2023-06-08 Thursday", $true, $true, $false, $false, $false, $true, 1, $false, "2023-06-09 Friday", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("63-12=51", $true, $true, $false, $false, $false, $true, 1, $false, "37-31=6", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("37-5=32", $true, $true, $false, $false, $false, $true, 1, $false, "42-22=20", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("25+37=62", $true, $true, $false, $false, $false, $true, 1, $false, "42+35=77", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("80-75=5", $true, $true, $false, $false, $false, $true, 1, $false, "99-55=44", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("95-24=71", $true, $true, $false, $false, $false, $true, 1, $false, "8+19=27", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("57-32=25", $true, $true, $false, $false, $false, $true, 1, $false, "70-32=38", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("45+5=50", $true, $true, $false, $false, $false, $true, 1, $false, "52+15=67", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("59-37=22", $true, $true, $false, $false, $false, $true, 1, $false, "16+46=62", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("1+1=2", $true, $true, $false, $false, $false, $true, 1, $false, "4+80=84", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("90-39=51", $true, $true, $false, $false, $false, $true, 1, $false, "52+18=70", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("31-9=22", $true, $true, $false, $false, $false, $true, 1, $false, "42+41=83", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("52+12=64", $true, $true, $false, $false, $false, $true, 1, $false, "11+86=97", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("24+40=64", $true, $true, $false, $false, $false, $true, 1, $false, "23-14=9", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("17+51=68", $true, $true, $false, $false, $false, $true, 1, $false, "69-33=36", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("2+10=12", $true, $true, $false, $false, $false, $true, 1, $false, "56+39=95", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("56+37=93", $true, $true, $false, $false, $false, $true, 1, $false, "27+13=40", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("43-2=41", $true, $true, $false, $false, $false, $true, 1, $false, "68+6=74", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("17+17=34", $true, $true, $false, $false, $false, $true, 1, $false, "74+19=93", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("92-72=20", $true, $true, $false, $false, $false, $true, 1, $false, "44-11=33", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("52+13=65", $true, $true, $false, $false, $false, $true, 1, $false, "62-21=41", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("40+48=88", $true, $true, $false, $false, $false, $true, 1, $false, "90-69=21", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("90-38=52", $true, $true, $false, $false, $false, $true, 1, $false, "6+60=66", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("38+40=78", $true, $true, $false, $false, $false, $true, 1, $false, "58-29=29", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("75-71=4", $true, $true, $false, $false, $false, $true, 1, $false, "76+11=87", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("49+14=63", $true, $true, $false, $false, $false, $true, 1, $false, "81-15=66", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("80-15=65", $true, $true, $false, $false, $false, $true, 1, $false, "90-6=84", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("44+44=88", $true, $true, $false, $false, $false, $true, 1, $false, "8+44=52", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("87-33=54", $true, $true, $false, $false, $false, $true, 1, $false, "42+12=54", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("30-15=15", $true, $true, $false, $false, $false, $true, 1, $false, "77+3=80", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("76-70=6", $true, $true, $false, $false, $false, $true, 1, $false, "77-35=42", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("81-21=60", $true, $true, $false, $false, $false, $true, 1, $false, "85-44=41", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("55+0=55", $true, $true, $false, $false, $false, $true, 1, $false, "66+28=94", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("28+67=95", $true, $true, $false, $false, $false, $true, 1, $false, "37-12=25", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("66-24=42", $true, $true, $false, $false, $false, $true, 1, $false, "3+3=6", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("66+14=80", $true, $true, $false, $false, $false, $true, 1, $false, "13+76=89", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("98-36=62", $true, $true, $false, $false, $false, $true, 1, $false, "56+23=79", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("25+29=54", $true, $true, $false, $false, $false, $true, 1, $false, "27+71=98", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("18+13=31", $true, $true, $false, $false, $false, $true, 1, $false, "53+25=78", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("14+62=76", $true, $true, $false, $false, $false, $true, 1, $false, "64-44=20", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("90-78=12", $true, $true, $false, $false, $false, $true, 1, $false, "69+28=97", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("79-45=34", $true, $true, $false, $false, $false, $true, 1, $false, "8+11=19", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("23+10=33", $true, $true, $false, $false, $false, $true, 1, $false, "78-65=13", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("97-18=79", $true, $true, $false, $false, $false, $true, 1, $false, "18-1=17", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("60-14=46", $true, $true, $false, $false, $false, $true, 1, $false, "68+25=93", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("0+37=37", $true, $true, $false, $false, $false, $true, 1, $false, "90-32=58", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("33+33=66", $true, $true, $false, $false, $false, $true, 1, $false, "49+16=65", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("18+45=63", $true, $true, $false, $false, $false, $true, 1, $false, "56-36=20", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("91-9=82", $true, $true, $false, $false, $false, $true, 1, $false, "78-15=63", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("66-17=49", $true, $true, $false, $false, $false, $true, 1, $false, "36-33=3", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("28+35=63", $true, $true, $false, $false, $false, $true, 1, $false, "86-20=66", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("28+13=41", $true, $true, $false, $false, $false, $true, 1, $false, "62-39=23", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("58-6=52", $true, $true, $false, $false, $false, $true, 1, $false, "12+84=96", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("4-2=2", $true, $true, $false, $false, $false, $true, 1, $false, "63-1=62", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("87-39=48", $true, $true, $false, $false, $false, $true, 1, $false, "69-42=27", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("56+41=97", $true, $true, $false, $false, $false, $true, 1, $false, "64-55=9", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("31-4=27", $true, $true, $false, $false, $false, $true, 1, $false, "81-66=15", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("60-17=43", $true, $true, $false, $false, $false, $true, 1, $false, "32+54=86", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("99-85=14", $true, $true, $false, $false, $false, $true, 1, $false, "92-87=5", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("97-53=44", $true, $true, $false, $false, $false, $true, 1, $false, "46-22=24", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("82-1=81", $true, $true, $false, $false, $false, $true, 1, $false, "50-26=24", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("15+59=74", $true, $true, $false, $false, $false, $true, 1, $false, "21+16=37", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("25+42=67", $true, $true, $false, $false, $false, $true, 1, $false, "26+29=55", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("80-70=10", $true, $true, $false, $false, $false, $true, 1, $false, "58-35=23", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("44+41=85", $true, $true, $false, $false, $false, $true, 1, $false, "45-22=23", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("92-58=34", $true, $true, $false, $false, $false, $true, 1, $false, "5+0=5", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("71-23=48", $true, $true, $false, $false, $false, $true, 1, $false, "79-61=18", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("88-13=75", $true, $true, $false, $false, $false, $true, 1, $false, "29+24=53", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("22+59=81", $true, $true, $false, $false, $false, $true, 1, $false, "97-38=59", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("19-5=14", $true, $true, $false, $false, $false, $true, 1, $false, "55+7=62", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("37-14=23", $true, $true, $false, $false, $false, $true, 1, $false, "69-19=50", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("67-64=3", $true, $true, $false, $false, $false, $true, 1, $false, "47+26=73", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("10+79=89", $true, $true, $false, $false, $false, $true, 1, $false, "79+15=94", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("95-89=6", $true, $true, $false, $false, $false, $true, 1, $false, "45-28=17", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("20-12=8", $true, $true, $false, $false, $false, $true, 1, $false, "20-2=18", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("23+15=38", $true, $true, $false, $false, $false, $true, 1, $false, "54-42=12", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("23+6=29", $true, $true, $false, $false, $false, $true, 1, $false, "18+50=68", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("43+9=52", $true, $true, $false, $false, $false, $true, 1, $false, "42-23=19", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("66-11=55", $true, $true, $false, $false, $false, $true, 1, $false, "28-13=15", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("61-53=8", $true, $true, $false, $false, $false, $true, 1, $false, "41-2=39", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("74+6=80", $true, $true, $false, $false, $false, $true, 1, $false, "70+26=96", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("96-62=34", $true, $true, $false, $false, $false, $true, 1, $false, "45-2=43", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("88+5=93", $true, $true, $false, $false, $false, $true, 1, $false, "96-2=94", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("77+4=81", $true, $true, $false, $false, $false, $true, 1, $false, "5+12=17", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("96-95=1", $true, $true, $false, $false, $false, $true, 1, $false, "75-24=51", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("59+4=63", $true, $true, $false, $false, $false, $true, 1, $false, "58-13=45", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("85-35=50", $true, $true, $false, $false, $false, $true, 1, $false, "81-28=53", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("44+1=45", $true, $true, $false, $false, $false, $true, 1, $false, "0+24=24", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("96+1=97", $true, $true, $false, $false, $false, $true, 1, $false, "29+68=97", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("18+41=59", $true, $true, $false, $false, $false, $true, 1, $false, "10+42=52", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("66-59=7", $true, $true, $false, $false, $false, $true, 1, $false, "0+21=21", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("68-25=43", $true, $true, $false, $false, $false, $true, 1, $false, "37+28=65", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("21+1=22", $true, $true, $false, $false, $false, $true, 1, $false, "16+11=27", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("98-53=45", $true, $true, $false, $false, $false, $true, 1, $false, "80-23=57", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("69+24=93", $true, $true, $false, $false, $false, $true, 1, $false, "51+20=71", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("86-81=5", $true, $true, $false, $false, $false, $true, 1, $false, "32+8=40", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("92-40=52", $true, $true, $false, $false, $false, $true, 1, $false, "28-11=17", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("17+40=57", $true, $true, $false, $false, $false, $true, 1, $false, "39-11=28", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("44+5=49", $true, $true, $false, $false, $false, $true, 1, $false, "35-32=3", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("63-51=12", $true, $true, $false, $false, $false, $true, 1, $false, "75-5=70", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("68+26=94", $true, $true, $false, $false, $false, $true, 1, $false, "87-82=5", 2) | Out-Null
$r = $d.Content
